$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 0) Locate the "Living in a Junkyard searching ..." paragraph up front;
#    everything below is scoped off of it so we never touch unrelated text.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$livingIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Living in a Junkyard*") {
        $livingIdx = $i
        break
    }
}

# ------------------------------------------------------------------
# 1) "Living in a Junkyard searching ..." -> "Living in a Junkyard, searching ..."
#    (insert a comma right after "Junkyard"), scoped to just that paragraph.
# ------------------------------------------------------------------
if ($livingIdx -gt 0) {
    $livingRange = $d.Paragraphs.Item($livingIdx).Range
    $found = $livingRange.Find.Execute("n a Junkyard ", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "n a Junkyard, ", 2)
}

# ------------------------------------------------------------------
# 2) Walk forward past two blank paragraphs; the 3rd blank paragraph
#    after it is the one that gets the new
#    "--Control is given to the player, ..." line.
# ------------------------------------------------------------------
if ($livingIdx -gt 0) {
    $blankSeen = 0
    $targetIdx = -1
    for ($i = $livingIdx + 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        # a "blank" paragraph contains nothing but the paragraph mark
        if ($t -eq "" -or $t -eq "`r") {
            $blankSeen = $blankSeen + 1
            if ($blankSeen -eq 3) {
                $targetIdx = $i
                break
            }
        }
    }

    if ($targetIdx -gt 0) {
        $targetPara = $d.Paragraphs.Item($targetIdx)
        $targetPara.Range.Text = "--Control is given to the player, player goes outside walks around and meets Ned the CEO of Garbage Trucks--"
    }
}
